$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: B3 becomes numeric 123 with the "left aligned" style (same style as B2/B4 originally)
$ws.Range("B3").Value = 123
$ws.Range("B3").HorizontalAlignment = -4131

# Row 4: A4 becomes "Admin", B4 becomes "admin123" with default style
$ws.Range("A4").Value = "Admin"
$ws.Range("B4").Value = "admin123"
$ws.Range("B4").Style = "Normal"

# Row 5: A5/B5 become the new credentials, no special style
$ws.Range("A5").Value = "qasmart.zee@bssuniversal.com"
$ws.Range("B5").Value = "Bss@2025-1"

# Row 6 gets removed entirely
$ws.Rows.Item(6).Delete() | Out-Null

# Update the shared string content used for the remaining cells
$ws.Range("B14").Select() | Out-Null
